$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column (D) cells whose new values look like
# plain numbers (e.g. "411.34"), so Excel keeps them as text strings
# instead of auto-converting them to numeric values, matching the source data.
$priceCells = @('D5', 'D6', 'D8', 'D10', 'D11', 'D12', 'D15', 'D16', 'D20', 'D23', 'D24', 'D25', 'D26', 'D27', 'D29', 'D30', 'D31', 'D34', 'D35', 'D38', 'D41', 'D43', 'D45', 'D46', 'D48', 'D49')
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range('D2').Value = '63.019.07'
$ws.Range('E2').Value = '  +6.68%  '
$ws.Range('D3').Value = '3.518.56'
$ws.Range('E3').Value = '  +6.17%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '411.34'
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('D6').Value = '129.46'
$ws.Range('E6').Value = '  +14.95%  '
$ws.Range('D7').Value = '3.509.43'
$ws.Range('E7').Value = '  +6.23%  '
$ws.Range('D8').Value = '0.597'
$ws.Range('E8').Value = '  +2.09%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').Value = '0.691'
$ws.Range('E10').Value = '  +9.24%  '
$ws.Range('D11').Value = '0.129'
$ws.Range('E11').Value = '  +31.43%  '
$ws.Range('D12').Value = '43.15'
$ws.Range('E12').Value = '  +7.79%  '
$ws.Range('D13').Value = '4.062.02'
$ws.Range('E13').Value = '  +5.48%  '
$ws.Range('E14').Value = '  -0.54%  '
$ws.Range('D15').Value = '8.74'
$ws.Range('E15').Value = '  +3.09%  '
$ws.Range('D16').Value = '20.22'
$ws.Range('E16').Value = '  +4.10%  '
$ws.Range('D17').Value = '3.472.84'
$ws.Range('E17').Value = '  +5.63%  '
$ws.Range('D18').Value = '62.979.49'
$ws.Range('E18').Value = '  +6.86%  '
$ws.Range('E19').Value = '  +0.97%  '
$ws.Range('D20').Value = '11.12'
$ws.Range('E20').Value = '  +3.97%  '
$ws.Range('E21').Value = '  +24.76%  '
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('D23').Value = '82.09'
$ws.Range('E23').Value = '  +9.14%  '
$ws.Range('D24').Value = '13.11'
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('D25').Value = '314.03'
$ws.Range('E25').Value = '  +3.42%  '
$ws.Range('D26').Value = '3.18'
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').Value = '30.60'
$ws.Range('E27').Value = '  +7.50%  '
$ws.Range('E28').Value = '  +4.11%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').Value = '0.182'
$ws.Range('E29').Value = '  +0.67%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').Value = '7.80'
$ws.Range('E30').Value = '  +3.05%  '
$ws.Range('D31').Value = '4.36'
$ws.Range('E31').Value = '  -2.51%  '
$ws.Range('E32').Value = '  +3.33%  '
$ws.Range('E33').Value = '  +5.52%  '
$ws.Range('D34').Value = '2.69'
$ws.Range('E34').Value = '  +28.23%  '
$ws.Range('D35').Value = '43.33'
$ws.Range('E35').Value = '  +9.48%  '
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('E37').Value = '  -4.06%  '
$ws.Range('D38').Value = '52.68'
$ws.Range('E38').Value = '  +1.26%  '
$ws.Range('E39').Value = '  +5.66%  '
$ws.Range('E40').Value = '  -0.26%  '
$ws.Range('D41').Value = '3.02'
$ws.Range('E41').Value = '  -2.89%  '
$ws.Range('E42').Value = '  +5.29%  '
$ws.Range('D43').Value = '138.26'
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('E44').Value = '  +2.78%  '
$ws.Range('D45').Value = '17.74'
$ws.Range('E45').Value = '  +4.88%  '
$ws.Range('D46').Value = '0.288'
$ws.Range('E46').Value = '  +1.62%  '
$ws.Range('E47').Value = '  +1.35%  '
$ws.Range('D48').Value = '2.25'
$ws.Range('E48').Value = '  -0.85%  '
$ws.Range('D49').Value = '22.42'
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('D50').Value = '2.222.35'
$ws.Range('E50').Value = '  +0.56%  '
$ws.Range('D51').Value = '3.859.75'
$ws.Range('E51').Value = '  +5.32%  '
